$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated price strings are purely numeric-looking (e.g. "5.59") and
# Excel would auto-coerce a plain Value assignment into a Number. The source
# data models these Price cells as text, so force Text formatting first on
# exactly the cells that need it, preserving the original text cell type.
$textPriceCells = @("D5", "D6", "D8", "D11", "D13", "D14", "D20", "D22", "D24", "D26", "D28", "D29", "D32", "D33", "D35", "D36", "D37", "D39", "D42", "D43", "D45", "D46", "D49", "D50", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.781.56"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "2.538.03"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "566.47"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "146.60"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").Value = "2.536.33"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "26.91"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "2.995.67"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "62.777.23"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "2.539.67"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "333.50"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "6.75"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "64.66"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").Value = "1.59"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "1.47"
$ws.Range("E28").Value = "  +10.68%  "
$ws.Range("D29").Value = "8.30"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "1.84"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "176.74"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("D35").Value = "404.23"
$ws.Range("E35").Value = "  +8.88%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "18.90"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.395"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "4.31"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "38.88"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").Value = "151.45"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "20.52"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "0.600"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("D50").Value = "18.20"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "1.76"
$ws.Range("E51").Value = "  +0.29%  "
